# "Generate Report for Handoff"
#
# Updates the localization-status report to reflect a fresh handoff run:
#   - Status changes from "Handed back: in sync with en-US" to
#     "Ready for handoff" everywhere it appears (Overview + per-locale sheets).
#   - The handoff timestamps are bumped to the new generation time.
#   - The (now much shorter) "Status" values no longer need as much column
#     room, so the Status-ish columns are narrowed to match.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-16 02:53:29"

# Columns E (zh-cn) and F (de-de) narrow from ~30 chars to ~17.2 chars.
$ov.Range("E1:F1").ColumnWidth = 16.29

# ---- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-08-16 02:53:25"
$zh.Range("C1").ColumnWidth = 16.29

# ---- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Ready for handoff"
$de.Range("H2").Value = "2016-08-16 02:53:29"
$de.Range("C1").ColumnWidth = 16.29
